$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(35).Copy()
$ws.Rows.Item(36).Insert()
Write-Host "done"
